$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices: A=1 B=2 C=3 D=4 E=5
# Using a leading apostrophe forces text entry (so numeric-looking strings like
# "345.19" stay text instead of becoming real numbers), and resetting the
# style back to "Normal" afterwards avoids Excel silently tagging the cell
# with a "quoted text" number format / style index.

function Set-Cell($row, $col, $val) {
    if ($val -eq $null) { return }
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

function Set-Row($row, $b, $c, $d, $e) {
    Set-Cell $row 2 $b
    Set-Cell $row 3 $c
    Set-Cell $row 4 $d
    Set-Cell $row 5 $e
}

Set-Row 2  $null $null "52.264.89" "  +1.58%  "
Set-Row 3  $null $null "2.796.78"  "  +1.78%  "
Set-Row 4  $null $null $null       "  +0.01%  "
Set-Row 5  $null $null "345.19"    "  +4.18%  "
Set-Row 6  $null $null "116.74"    "  +0.86%  "
Set-Row 7  $null $null "0.554"     "  +4.15%  "
Set-Row 8  $null $null $null       "  -0.02%  "
Set-Row 9  $null $null "0.586"     "  +2.75%  "
Set-Row 10 $null $null "43.18"     "  +3.66%  "
Set-Row 11 $null $null "0.0859"    "  +3.28%  "
Set-Row 12 $null $null "20.20"     "  -0.91%  "
Set-Row 13 $null $null $null       "  +1.83%  "
Set-Row 14 $null $null "7.79"      "  +1.39%  "
Set-Row 15 $null $null "3.235.95"  "  +1.98%  "
Set-Row 16 $null $null "2.799.09"  "  +2.14%  "
Set-Row 17 $null $null "0.893"     "  +0.69%  "
Set-Row 18 $null $null "52.161.91" "  +1.49%  "
Set-Row 19 $null $null "3.21"      "  +6.35%  "
Set-Row 20 $null $null "7.10"      "  +3.60%  "
Set-Row 21 $null $null "13.45"     "  -1.45%  "
Set-Row 22 $null $null "0.0₃0984"  "  +2.18%  "
Set-Row 23 $null $null "70.33"     "  -0.37%  "
Set-Row 24 $null $null $null       "  -6.35%  "
Set-Row 25 $null $null "2.76"      "  +6.34%  "
Set-Row 26 $null $null "26.71"     "  -0.74%  "
Set-Row 27 $null $null $null       "  -0.01%  "
Set-Row 28 $null $null "10.30"     "  -0.63%  "
Set-Row 29 $null $null $null       "  +0.51%  "
Set-Row 30 $null $null $null       "  -1.15%  "
Set-Row 31 $null $null "35.07"     "  -1.47%  "
Set-Row 32 $null $null "50.26"     "  +0.06%  "
Set-Row 33 $null $null "5.75"      "  +2.31%  "

Set-Row 34 "VeChain"         "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"            "0.0416" "  +18.27%  "
Set-Row 35 "Hedera"          "https://coinranking.com/coin/jad286TjB+hedera-hbar"                 "0.0826" "  -0.24%  "
Set-Row 36 "ARBITRUM"        "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"                "2.12"   "  +0.66%  "
Set-Row 37 "FirstDigitalUSD" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"       "1.00"   "  -0.15%  "

Set-Row 38 $null $null "19.07" "  -2.21%  "
Set-Row 39 $null $null "4.98"  "  -1.11%  "
Set-Row 40 $null $null "3.22"  "  -0.03%  "
Set-Row 41 $null $null "2.71"  "  +21.76%  "
Set-Row 42 $null $null "23.61" "  -0.35%  "

Set-Row 43 "Stellar" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" "0.116"  "  +2.43%  "
Set-Row 44 "Monero"  "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" "127.52" "  -1.63%  "

Set-Row 45 $null $null $null  "  +1.68%  "
Set-Row 46 $null $null "3.36" "  -1.67%  "
Set-Row 47 $null $null "2.080.86" "  -1.28%  "
Set-Row 48 $null $null "2.30" "  +3.37%  "

Set-Row 49 "SEI"       "https://coinranking.com/coin/8nxCqs-uj+sei-sei"      "0.963" "  +16.07%  "
Set-Row 50 "THORChain" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune" "5.56"  "  +0.52%  "

Set-Row 51 $null $null "8.98" "  -1.42%  "
